# Insert a new daily price record for "Acelga" (Vega Modelo de Temuco) as
# row 309, pushing the existing rows 309-398 down to 310-399.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 309 - everything below shifts
# down by one (Excel also extends the used range / dimension to R399).
$ws.Rows("309:309").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A309").Value = 10
$ws.Range("B309").Value = "Vega Modelo de Temuco"
$ws.Range("C309").Value = "La Araucanía"
$ws.Range("D309").Value = 44841
$ws.Range("E309").Value = 9
$ws.Range("F309").Value = 100112009
$ws.Range("G309").Value = "Acelga"
$ws.Range("H309").Value = "Sin especificar"
$ws.Range("I309").Value = "Primera"
$ws.Range("J309").Value = 70
$ws.Range("K309").Value = 8000
$ws.Range("L309").Value = 8000
$ws.Range("M309").Value = 8000
$ws.Range("N309").Value = "$/docena de atados (12 kilos)"
$ws.Range("O309").Value = "Provincia de Cautín"
$ws.Range("P309").Value = 667
$ws.Range("Q309").Value = 12
$ws.Range("R309").Value = "Hortaliza"
